$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header for new column K
$ws.Range("K1").Value = "intervention_type"

# Copy the header style (bold, border, centered) from an existing header cell (J1) to K1
$ws.Range("J1").Copy()
$ws.Range("K1").PasteSpecial(-4122) # xlPasteFormats
$excel.CutCopyMode = $false

# Values for the new "intervention_type" column, rows 2-16
$values = @(
    "PROCEDURE",
    "DRUG",
    "DEVICE",
    "DIAGNOSTIC_TEST",
    "OTHER",
    "OTHER",
    "OTHER",
    "BEHAVIORAL",
    "DIAGNOSTIC_TEST",
    "OTHER",
    "COMBINATION_PRODUCT",
    "DIAGNOSTIC_TEST",
    "OTHER",
    "OTHER",
    "DEVICE"
)

for ($i = 0; $i -lt $values.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 11).Value = $values[$i]
}
